$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("protocol")
$ws2 = $wb.Worksheets.Item("protocoltestcasedetails")

# ---- protocoltestcasedetails: drop the two trailing blank rows (40,41) ----
$ws2.Range("A40:C41").Delete(-4162)

# ---- fix the stray "header" style that leaked onto blank row 35 ----
$ws2.Range("A34:C34").Copy()
$ws2.Range("A35:C35").PasteSpecial(-4122)

# ---- new column D: "execute" flag ----
# Enter the data in the same order the original authoring session did
# (N first down the column, Y on the one exception row, then the header)
# so the shared-string table comes out in the same order.
for ($r = 2; $r -le 27; $r++) {
    if ($r -eq 23) {
        $ws2.Cells.Item($r, 4).Value = "Y"
    } else {
        $ws2.Cells.Item($r, 4).Value = "N"
    }
}
$ws2.Range("D1").Value = "execute"

# Formatting: D2:D39 take the bordered/unlocked look of column C, centered
$ws2.Range("C2").Copy()
$ws2.Range("D2:D39").PasteSpecial(-4122)
$ws2.Range("D2:D39").HorizontalAlignment = -4108

# D1 gets the bold header look shared by A1:C1
$ws2.Range("C1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)

# Column D width
$ws2.Columns("D").ColumnWidth = 7.25

# Dropdown validation list restricting D2:D27 to Y/N
$ws2.Range("D2:D27").Validation.Add(3, 1, 1, """Y,N""")

# ---- sheet views / selections ----
$ws1.Range("B12").Select()
$ws2.Select()
$ws2.Range("C28").Select()
$excel.ActiveWindow.ScrollRow = 16

Write-Host "edit applied"
